$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3F")

# Fix C5: was inline string "05", now numeric 5
$ws.Range("C5").Value = 5

# Add new row 6 for the new submission
$ws.Range("A6").Value = "2026-02-08 22:15:44"
$ws.Range("B6").Value = "Usman Muhammad Gubio"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "05"
$ws.Range("D6").Value = 7
